$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark (it currently sits on the
#    last, empty paragraph of the document). It gets re-created further
#    below on the newly inserted paragraph, so delete the old one first
#    to avoid any name collision with the new one we are about to add.
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# 2) Insert a brand-new paragraph right after the first paragraph
#    ("Entity Framework:") containing the bold "aaaaaaaaaaa" run and a
#    fresh "_GoBack" bookmark. We target an insertion point just before
#    the last character of paragraph 1's text (not the paragraph-end
#    boundary itself) so the existing paragraph/run is left completely
#    untouched and a clean new paragraph is spliced in right after it.
$firstPara = $d.Paragraphs(1)
$insertAt = $firstPara.Range.End - 1
$insertionRange = $d.Range($insertAt, $insertAt)

$newParagraphXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr>' +
  '<w:rPr>' +
  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/>' +
  '<w:sz w:val="24"/>' +
  '<w:szCs w:val="24"/>' +
  '</w:rPr>' +
  '</w:pPr>' +
  '<w:r>' +
  '<w:rPr>' +
  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:b/>' +
  '<w:sz w:val="24"/>' +
  '<w:szCs w:val="24"/>' +
  '</w:rPr>' +
  '<w:t>aaaaaaaaaaa</w:t>' +
  '</w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$insertionRange.InsertXML($newParagraphXml)
